$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(90133, "Dra. Beatriz da Mata", "Vendas", "Viagem de negócios", 2, 45102, 2782.67),
    @(80059, "Lucas Cardoso", "P&D", "Problemas pessoais", 3, 45104, 11899.41),
    @(74569, "Luiz Otávio Monteiro", "Recursos Humanos", "Viagem de negócios", 1, 45106, 12423.22),
    @(63240, "Clarice Correia", "Engenharia", "Viagem de negócios", 2, 45095, 8390.709999999999),
    @(75449, "Mariana Porto", "Jurídico", "Consulta médica", 3, 45086, 2772.57),
    @(24348, "Dr. Luiz Gustavo Almeida", "Jurídico", "Viagem de negócios", 4, 45087, 11681.88),
    @(68655, "Danilo Caldeira", "Recursos Humanos", "Problemas pessoais", 6, 45092, 8639.77),
    @(90731, "Amanda Melo", "Marketing", "Doença", 2, 45103, 9703.82),
    @(25362, "Maria Luiza Freitas", "Atendimento ao Cliente", "Doença", 7, 45086, 3196.1),
    @(35913, "Francisco Moraes", "Vendas", "Problemas pessoais", 4, 45104, 11202.75)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $rowData = $data[$i]
    $ws.Cells.Item($row, 1).Value = $rowData[0]
    $ws.Cells.Item($row, 2).Value = $rowData[1]
    $ws.Cells.Item($row, 3).Value = $rowData[2]
    $ws.Cells.Item($row, 4).Value = $rowData[3]
    $ws.Cells.Item($row, 5).Value = $rowData[4]
    $ws.Cells.Item($row, 6).Value = $rowData[5]
    $ws.Cells.Item($row, 7).Value = $rowData[6]
}
